$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Estimated")
$ws.Range("B2").Value = "Tue Feb 06 16:19:47 EST 2024"
$ws.Range("B3").Value = "Tue Feb 06 16:20:45 EST 2024"
$ws.Range("B4").Value = "Tue Feb 06 16:21:39 EST 2024"
$ws.Range("B5").Value = "Tue Feb 06 16:22:33 EST 2024"
$ws.Range("B6").Value = "Tue Feb 06 16:23:30 EST 2024"
$ws.Range("B7").Value = "Tue Feb 06 16:24:24 EST 2024"

$ws = $wb.Worksheets.Item("Existing")
$ws.Range("B2").Value = "Tue Feb 06 16:25:17 EST 2024"
$ws.Range("B3").Value = "Tue Feb 06 16:26:12 EST 2024"
$ws.Range("B4").Value = "Tue Feb 06 16:27:06 EST 2024"
$ws.Range("B5").Value = "Tue Feb 06 16:27:59 EST 2024"
$ws.Range("B6").Value = "Tue Feb 06 16:28:52 EST 2024"
$ws.Range("B7").Value = "Tue Feb 06 16:29:45 EST 2024"
$ws.Range("B8").Value = "Tue Feb 06 16:30:39 EST 2024"
$ws.Range("B9").Value = "Tue Feb 06 16:31:33 EST 2024"
$ws.Range("B10").Value = "Tue Feb 06 16:32:26 EST 2024"
$ws.Range("B11").Value = "Tue Feb 06 16:33:23 EST 2024"
$ws.Range("B12").Value = "Tue Feb 06 16:34:22 EST 2024"
$ws.Range("B13").Value = "Tue Feb 06 16:35:16 EST 2024"
$ws.Range("B14").Value = "Tue Feb 06 16:36:12 EST 2024"
$ws.Range("B15").Value = "Tue Feb 06 16:37:08 EST 2024"
$ws.Range("B16").Value = "Tue Feb 06 16:38:02 EST 2024"
$ws.Range("B17").Value = "Tue Feb 06 16:38:56 EST 2024"
$ws.Range("B18").Value = "Tue Feb 06 16:39:48 EST 2024"
$ws.Range("B19").Value = "Tue Feb 06 16:40:41 EST 2024"

$ws = $wb.Worksheets.Item("NewTaxReturn")
$ws.Range("B2").Value = "Tue Feb 06 16:41:36 EST 2024"
$ws.Range("B3").Value = "Tue Feb 06 16:42:28 EST 2024"
$ws.Range("B4").Value = "Tue Feb 06 16:43:20 EST 2024"
$ws.Range("B5").Value = "Tue Feb 06 16:44:13 EST 2024"
$ws.Range("B6").Value = "Tue Feb 06 16:45:08 EST 2024"
$ws.Range("B7").Value = "Tue Feb 06 16:46:00 EST 2024"
$ws.Range("B8").Value = "Tue Feb 06 16:46:54 EST 2024"
$ws.Range("B9").Value = "Tue Feb 06 16:47:46 EST 2024"
$ws.Range("B10").Value = "Tue Feb 06 16:48:38 EST 2024"
$ws.Range("B11").Value = "Tue Feb 06 16:49:31 EST 2024"
$ws.Range("B12").Value = "Tue Feb 06 16:50:23 EST 2024"
$ws.Range("B13").Value = "Tue Feb 06 16:51:16 EST 2024"
$ws.Range("B14").Value = "Tue Feb 06 16:52:08 EST 2024"
$ws.Range("B15").Value = "Tue Feb 06 16:53:01 EST 2024"
$ws.Range("B16").Value = "Tue Feb 06 16:53:53 EST 2024"
$ws.Range("B17").Value = "Tue Feb 06 16:54:45 EST 2024"
$ws.Range("B18").Value = "Tue Feb 06 16:55:38 EST 2024"
$ws.Range("B19").Value = "Tue Feb 06 16:56:31 EST 2024"
$ws.Range("B20").Value = "Tue Feb 06 16:57:23 EST 2024"
$ws.Range("B21").Value = "Tue Feb 06 16:58:18 EST 2024"
$ws.Range("B22").Value = "Tue Feb 06 16:59:10 EST 2024"
$ws.Range("B23").Value = "Tue Feb 06 17:00:02 EST 2024"
$ws.Range("B24").Value = "Tue Feb 06 17:00:56 EST 2024"
$ws.Range("B25").Value = "Tue Feb 06 17:01:48 EST 2024"
$ws.Range("B26").Value = "Tue Feb 06 17:02:40 EST 2024"
$ws.Range("B27").Value = "Tue Feb 06 17:03:32 EST 2024"
$ws.Range("B28").Value = "Tue Feb 06 17:04:24 EST 2024"
$ws.Range("B29").Value = "Tue Feb 06 17:05:16 EST 2024"
$ws.Range("B30").Value = "Tue Feb 06 17:06:08 EST 2024"
$ws.Range("B31").Value = "Tue Feb 06 17:07:01 EST 2024"
$ws.Range("B32").Value = "Tue Feb 06 17:07:54 EST 2024"
$ws.Range("B33").Value = "Tue Feb 06 17:08:47 EST 2024"
$ws.Range("B34").Value = "Tue Feb 06 17:09:39 EST 2024"
$ws.Range("B35").Value = "Tue Feb 06 17:10:32 EST 2024"
$ws.Range("B36").Value = "Tue Feb 06 17:11:26 EST 2024"
$ws.Range("B37").Value = "Tue Feb 06 17:12:18 EST 2024"
$ws.Range("B38").Value = "Tue Feb 06 17:13:11 EST 2024"
$ws.Range("B39").Value = "Tue Feb 06 17:14:03 EST 2024"
$ws.Range("B40").Value = "Tue Feb 06 17:14:56 EST 2024"
$ws.Range("B41").Value = "Tue Feb 06 17:15:48 EST 2024"
$ws.Range("B42").Value = "Tue Feb 06 17:16:40 EST 2024"
$ws.Range("B43").Value = "Tue Feb 06 17:17:34 EST 2024"
$ws.Range("B44").Value = "Tue Feb 06 17:18:27 EST 2024"
$ws.Range("B45").Value = "Tue Feb 06 17:19:19 EST 2024"
$ws.Range("B46").Value = "Tue Feb 06 17:20:12 EST 2024"
$ws.Range("B47").Value = "Tue Feb 06 17:21:04 EST 2024"
$ws.Range("B48").Value = "Tue Feb 06 17:21:57 EST 2024"
$ws.Range("B49").Value = "Tue Feb 06 17:22:50 EST 2024"
$ws.Range("B50").Value = "Tue Feb 06 17:23:42 EST 2024"
$ws.Range("B51").Value = "Tue Feb 06 17:24:35 EST 2024"
$ws.Range("B52").Value = "Tue Feb 06 17:25:27 EST 2024"

$ws = $wb.Worksheets.Item("Personal_IND")
$ws.Range("B2").Value = "Tue Feb 06 17:28:03 EST 2024"
$ws.Range("B4").Value = "Tue Feb 06 17:28:53 EST 2024"
$ws.Range("B5").Value = "Tue Feb 06 17:29:43 EST 2024"
$ws.Range("B6").Value = "Tue Feb 06 17:30:35 EST 2024"
$ws.Range("B7").Value = "Tue Feb 06 17:31:32 EST 2024"
$ws.Range("B8").Value = "Tue Feb 06 17:32:24 EST 2024"
$ws.Range("B9").Value = "Tue Feb 06 17:33:14 EST 2024"

$ws = $wb.Worksheets.Item("Personal_JNT")
$ws.Range("B2").Value = "Tue Feb 06 17:34:05 EST 2024"
$ws.Range("B4").Value = "Tue Feb 06 17:35:00 EST 2024"
$ws.Range("B5").Value = "Tue Feb 06 17:35:56 EST 2024"
$ws.Range("B6").Value = "Tue Feb 06 17:36:54 EST 2024"

$ws = $wb.Worksheets.Item("Personal_EL")
$ws.Range("B2").Value = "Tue Feb 06 17:26:20 EST 2024"
$ws.Range("B3").Value = "Tue Feb 06 17:27:11 EST 2024"
